$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.990.47"
$ws.Range("E2").Value = "  -0.08%  "
$ws.Range("D3").Value = "1.919.32"
$ws.Range("E3").Value = "  +0.51%  "
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "325.52"
$ws.Range("E5").Value = "  +0.15%  "
$ws.Range("E6").Value = "  +0.16%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4587"
$ws.Range("E7").Value = "  -0.39%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3811"
$ws.Range("E8").Value = "  -0.62%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07749"
$ws.Range("E9").Value = "  -0.08%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9774"
$ws.Range("E10").Value = "  -0.94%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "22.60"
$ws.Range("E11").Value = "  +2.18%  "
$ws.Range("D12").Value = "1.944.34"
$ws.Range("E12").Value = "  +1.35%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.706"
$ws.Range("E13").Value = "  -0.11%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.947"
$ws.Range("E14").Value = "  -0.86%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.07002"
$ws.Range("E15").Value = "  -1.10%  "
$ws.Range("E16").Value = "  +0.24%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "84.51"
$ws.Range("E17").Value = "  +0.44%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000009497"
$ws.Range("E18").Value = "  -0.88%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "16.66"
$ws.Range("E19").Value = "  -0.79%  "
$ws.Range("E20").Value = "  +0.08%  "
$ws.Range("D21").Value = "29.005.48"
$ws.Range("E21").Value = "  -0.07%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.349"
$ws.Range("E22").Value = "  +0.12%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.04"
$ws.Range("E23").Value = "  +0.52%  "
$ws.Range("D24").Value = "2.168.27"
$ws.Range("E24").Value = "  +1.03%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.070"
$ws.Range("E25").Value = "  -0.42%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "157.66"
$ws.Range("E26").Value = "  +0.67%  "
$ws.Range("E27").Value = "  -0.99%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.611"
$ws.Range("E28").Value = "  -0.23%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "117.63"
$ws.Range("E29").Value = "  -0.43%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.833"
$ws.Range("E30").Value = "  -0.16%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.09323"
$ws.Range("E31").Value = "  +0.41%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.8598"
$ws.Range("E32").Value = "  -0.66%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.086"
$ws.Range("E33").Value = "  -0.70%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.240"
$ws.Range("E34").Value = "  -1.34%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.018"
$ws.Range("E35").Value = "  -0.04%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.05680"
$ws.Range("E36").Value = "  -1.22%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.153"
$ws.Range("E37").Value = "  -0.08%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.005"
$ws.Range("E38").Value = "  +0.30%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.161"
$ws.Range("E39").Value = "  +15.51%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.02039"
$ws.Range("E40").Value = "  -0.76%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "7.418"
$ws.Range("E41").Value = "  -1.17%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5488"
$ws.Range("E42").Value = "  -1.13%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1754"
$ws.Range("E43").Value = "  -0.49%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "9.361"
$ws.Range("E44").Value = "  +0.44%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.000002833"
$ws.Range("E45").Value = "  +7.04%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.181"
$ws.Range("E46").Value = "  +3.01%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5177"
$ws.Range("E47").Value = "  -1.00%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "11.14"
$ws.Range("E48").Value = "  -1.37%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.06918"
$ws.Range("E49").Value = "  +1.21%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "110.21"
$ws.Range("E50").Value = "  -1.76%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.758"
$ws.Range("E51").Value = "  -1.33%  "
